$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Nakul"
$ws.Cells.Item(2, 2).Value = 1919.7
$ws.Cells.Item(3, 1).Value = "Eric"
$ws.Cells.Item(3, 2).Value = 1850.2
$ws.Cells.Item(4, 1).Value = "Alexis"
$ws.Cells.Item(4, 2).Value = 1761.4
$ws.Cells.Item(5, 1).Value = "Sean B"
$ws.Cells.Item(5, 2).Value = 1760.4
$ws.Cells.Item(6, 1).Value = "Hassan"
$ws.Cells.Item(6, 2).Value = 1744.4
$ws.Cells.Item(7, 1).Value = "Saravanan"
$ws.Cells.Item(7, 2).Value = 1731.5
$ws.Cells.Item(8, 1).Value = "Paulo"
$ws.Cells.Item(8, 2).Value = 1712.3
$ws.Cells.Item(9, 1).Value = "Anil"
$ws.Cells.Item(9, 2).Value = 1685.6
$ws.Cells.Item(10, 1).Value = "Xi"
$ws.Cells.Item(10, 2).Value = 1684.9
$ws.Cells.Item(11, 1).Value = "Peter"
$ws.Cells.Item(11, 2).Value = 1667.8
$ws.Cells.Item(12, 1).Value = "Ayman"
$ws.Cells.Item(12, 2).Value = 1661.4
$ws.Cells.Item(13, 1).Value = "Charlie S"
$ws.Cells.Item(13, 2).Value = 1646.5
$ws.Cells.Item(14, 1).Value = "Victor"
$ws.Cells.Item(14, 2).Value = 1620
$ws.Cells.Item(15, 1).Value = "Mark"
$ws.Cells.Item(15, 2).Value = 1618.8
$ws.Cells.Item(16, 1).Value = "Zhengnan"
$ws.Cells.Item(16, 2).Value = 1616
$ws.Cells.Item(17, 1).Value = "Eugene"
$ws.Cells.Item(17, 2).Value = 1597.1
$ws.Cells.Item(18, 1).Value = "Luis"
$ws.Cells.Item(18, 2).Value = 1568.8
$ws.Cells.Item(19, 1).Value = "Karla"
$ws.Cells.Item(19, 2).Value = 1555
$ws.Cells.Item(20, 1).Value = "Sean H"
$ws.Cells.Item(20, 2).Value = 1552.6
$ws.Cells.Item(21, 1).Value = "Jameel"
$ws.Cells.Item(21, 2).Value = 1533.5
$ws.Cells.Item(22, 1).Value = "Dylan"
$ws.Cells.Item(22, 2).Value = 1519.2
$ws.Cells.Item(23, 1).Value = "Yevhen"
$ws.Cells.Item(23, 2).Value = 1513.2
$ws.Cells.Item(24, 1).Value = "Madi"
$ws.Cells.Item(24, 2).Value = 1512.5
$ws.Cells.Item(25, 1).Value = "Fernando"
$ws.Cells.Item(25, 2).Value = 1511.7
$ws.Cells.Item(26, 1).Value = "Sadeed"
$ws.Cells.Item(26, 2).Value = 1502.4
$ws.Cells.Item(27, 1).Value = "Octavio"
$ws.Cells.Item(27, 2).Value = 1495.4
$ws.Cells.Item(28, 1).Value = "Abdurauf"
$ws.Cells.Item(28, 2).Value = 1450.8
$ws.Cells.Item(29, 1).Value = "Jofrey"
$ws.Cells.Item(29, 2).Value = 1442.4
$ws.Cells.Item(30, 1).Value = "Amanat"
$ws.Cells.Item(30, 2).Value = 1438
$ws.Cells.Item(31, 1).Value = "Rawan"
$ws.Cells.Item(31, 2).Value = 1425.8
$ws.Cells.Item(32, 1).Value = "Carlos"
$ws.Cells.Item(32, 2).Value = 1417.1
$ws.Cells.Item(33, 1).Value = "Omar"
$ws.Cells.Item(33, 2).Value = 1407.5
$ws.Cells.Item(34, 1).Value = "Faruk"
$ws.Cells.Item(34, 2).Value = 1395
$ws.Cells.Item(35, 1).Value = "Mustafa"
$ws.Cells.Item(35, 2).Value = 1392.4
$ws.Cells.Item(36, 1).Value = "Abdulmajeed"
$ws.Cells.Item(36, 2).Value = 1370.5
$ws.Cells.Item(37, 1).Value = "Haytham"
$ws.Cells.Item(37, 2).Value = 1359.2
$ws.Cells.Item(38, 1).Value = "Lucas"
$ws.Cells.Item(38, 2).Value = 1335.4
$ws.Cells.Item(39, 1).Value = "Huawen"
$ws.Cells.Item(39, 2).Value = 1222.5
$ws.Cells.Item(40, 1).Value = "Xingzhu"
$ws.Cells.Item(40, 2).Value = 1214.9
$ws.Cells.Item(41, 1).Value = "Juris"
$ws.Cells.Item(41, 2).Value = 1211.4
$ws.Cells.Item(42, 1).Value = "Hashim"
$ws.Cells.Item(42, 2).Value = 1210.7
